# Append a second "After merging models" results table below the existing
# Precision/Recall table, mirroring the original layout (rows 1-10) in
# rows 14-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section title
$ws.Range("A14").Value = "After merging models"

# Header row (copy of row 1)
$ws.Range("A15").Value = "File Name"
$ws.Range("B15").Value = "Upto Frame Count"
$ws.Range("C15").Value = "TP"
$ws.Range("D15").Value = "FP"
$ws.Range("E15").Value = "TN"
$ws.Range("F15").Value = "FN"
$ws.Range("B15:F15").HorizontalAlignment = -4108

# Data rows (updated metrics after merging models)
$ws.Range("A16").Value = "combined_short.mp4"
$ws.Range("B16").Value = 9
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 0

$ws.Range("A17").Value = "normal_people_video.mp4"
$ws.Range("B17").Value = 28
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 109
$ws.Range("F17").Value = 0

$ws.Range("A18").Value = "worker_wrong_helmet_pos.mp4"
$ws.Range("B18").Value = 56
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 56
$ws.Range("F18").Value = 0

$ws.Range("A19").Value = "worker_without_helmet - Trim.mp4"
$ws.Range("B19").Value = 80
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 28
$ws.Range("E19").Value = 33
$ws.Range("F19").Value = 0

$ws.Range("A20").Value = "worker_with_helmet_multiple.mp4"
$ws.Range("B20").Value = 120
$ws.Range("C20").Value = 249
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 54

$ws.Range("B16:F20").HorizontalAlignment = -4108

# Totals row
$ws.Range("B21").Formula = "=SUM(B16:B20)"
$ws.Range("C21").Formula = "=SUM(C16:C20)"
$ws.Range("D21").Formula = "=SUM(D16:D20)"
$ws.Range("E21").Formula = "=SUM(E16:E20)"
$ws.Range("F21").Formula = "=SUM(F16:F20)"
$ws.Range("B21:F21").HorizontalAlignment = -4108

# Precision / Recall summary
$ws.Range("C23").Value = "Precision"
$ws.Range("D23").Formula = "=C21/SUM(C21+D21)"
$ws.Range("C24").Value = "Recall"
$ws.Range("D24").Formula = "=C21/SUM(C21+F21)"

# Update selection to match the saved file (selection at J16)
$ws.Range("J16").Select()

$wb.Save()
